$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the (empty) trailing column so the sheet's column count stays capped
# at the worksheet's real maximum once we shift everything right by one.
$ws.Columns("XFD:XFD").Delete()

# Insert a new blank column at B. Because columns A and B already share the
# same width/style, the new column merges into that existing span instead
# of needing a brand-new (unformatted) column definition.
$ws.Columns("B:B").Insert()

# The original column A content is still sitting in column A (now shifted
# logically to the "second" data column) - move it into the new column B.
$ws.Range("A1:A5").Copy()
$ws.Range("B1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Populate the new "TestCase" column in A.
$ws.Range("A1").Value = "TestCase"
$ws.Range("A2").Value = "FULL"
$ws.Range("A3").Value = "MANDATORY"
$ws.Range("A4").Value = "ADDITIONAL"
$ws.Range("A5").Value = "INVALID"
